# Adding the Shortcut-Hotkey Extension
# Adds a "2. Woche" section (row 10, merged A10:F10, same look as the
# "1. Woche" header in row 2) followed by two new task rows (11 and 12:
# "ShortcutDialog" and "Folder-Observer"), extends the shared E-column
# formula range down through the new rows, and pads a handful of empty
# rows below (13-21) so the sheet mirrors the rest of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 10: new "2. Woche" section header -- same formatting as A2:F2
# ---------------------------------------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("A10").Value = "2. Woche"
$ws.Range("A10:F10").Merge()
$ws.Rows.Item(10).RowHeight = 24.95

# ---------------------------------------------------------------
# Row 11: "ShortcutDialog" task
# ---------------------------------------------------------------
$ws.Range("A11").Value = "ShortcutDialog"
$ws.Range("B11").Value = 0.13541666666666666
$ws.Range("D11").Value = 0.0069444444444444441
$ws.Range("D11").NumberFormat = "h:mm"
$ws.Range("E11").Formula = '=IF(C11="",B11,C11) - D11'
$ws.Range("E11").NumberFormat = "hh:mm;@"
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("F11").Value = "David Eilmsteiner"
$ws.Rows.Item(11).RowHeight = 24.95

# ---------------------------------------------------------------
# Row 12: "Folder-Observer" task
# ---------------------------------------------------------------
$ws.Range("A12").Value = "Folder-Observer"
$ws.Range("B12").Value = 0.041666666666666664
$ws.Range("D12").Value = 0
$ws.Range("D12").NumberFormat = "h:mm"
$ws.Range("E12").Formula = '=IF(C12="",B12,C12) - D12'
$ws.Range("E12").NumberFormat = "hh:mm;@"
$ws.Range("E12").VerticalAlignment = -4108
$ws.Range("F12").Value = "Martin Niederl"
$ws.Rows.Item(12).RowHeight = 24.95

# ---------------------------------------------------------------
# Extend the shared "Remain" formula down through the new rows
# ---------------------------------------------------------------
$ws.Range("E4:E12").FormulaR1C1 = '=IF(RC3="",RC2,RC3) - RC4'

# ---------------------------------------------------------------
# Rows 13-21: trailing blank rows, matching existing look (E column
# keeps the time-remain formatting used throughout the sheet)
# ---------------------------------------------------------------
for ($r = 13; $r -le 21; $r++) {
    $ws.Range("E$r").NumberFormat = "hh:mm;@"
    $ws.Range("E$r").VerticalAlignment = -4108
    $ws.Rows.Item($r).RowHeight = 24.95
}

# ---------------------------------------------------------------
# Selection, matching what the author had selected on save
# ---------------------------------------------------------------
$ws.Range("G15").Select()

$wb.Saved = $false
